# Generate Report for Handoff
# Replaces the old handback/generate run (13d3d1d7-...) with the new one
# (7117ddf5-...), refreshes the associated hash-qualified .xlf file names,
# and bumps the handoff timestamps, across the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "13d3d1d7-1c30-4e30-b45c-2436e08b1062"
$newGuid = "7117ddf5-94d8-4ef0-9965-8f57d18f4f94"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldMdPath = "e2e\$oldGuid.md"
$newMdPath = "e2e\$newGuid.md"

$oldZhXlf = "$oldGuid.917df349decb8f55316f3e37b40ab0ce3cbb8a64.zh-cn.xlf"
$newZhXlf = "$newGuid.9f6e45ecc91e2356bd97aced1df829651e4e8364.zh-cn.xlf"

$oldDeXlf = "$oldGuid.917df349decb8f55316f3e37b40ab0ce3cbb8a64.de-de.xlf"
$newDeXlf = "$newGuid.9f6e45ecc91e2356bd97aced1df829651e4e8364.de-de.xlf"

# Hyperlink target (unchanged address, only the displayed text changes)
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc15a0a3827246b81ac16760456ec22d0a240ab5/e2e/$oldGuid.md"

function Update-Hyperlink($ws, $cellRef, $displayText) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellRef), $hyperlinkAddress, "", "", $displayText) | Out-Null
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = "2016-09-05 17:07:16"
Update-Hyperlink $wsOverview "B2" $newMdPath

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-09-05 17:07:11"
Update-Hyperlink $wsZhCn "A2" $newMd

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-09-05 17:07:16"
Update-Hyperlink $wsDeDe "A2" $newMd
